# Rename the "SwateTemplateMetadata" sheet to "isa_template" and update
# that sheet's active selection to C17 (matching the saved view state).

$wb = $excel.ActiveWorkbook

$metaSheet = $wb.Worksheets.Item("SwateTemplateMetadata")
$metaSheet.Name = "isa_template"

$metaSheet.Activate()
$metaSheet.Range("C17").Select()
